$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first two rows (Burnley fixture, West Ham away fixture)
$ws.Range("A1:B2").Delete()

# Insert a new row for the Club Brugge Champions League fixture after the
# Crystal Palace row (row 1 after deletion)
$ws.Range("A2:B2").Insert()
$ws.Range("A2").Value = "Manchester City v Club Brugge"
$ws.Range("B2").Value = "03/11/2021 20:00 | UEFA Champions League"

# Update the Wolverhampton Wanderers kickoff time (now row 5 after deletion+insert)
$ws.Range("B5").Value = "11/12/2021 12:30 | Premier League"

# Update the Leeds United fixture date (now row 6 after deletion+insert)
$ws.Range("B6").Value = "14/12/2021 20:00 | Premier League"
